$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "30.326.57"
$ws.Range("E2").Value = "  -1.10%  "

# Row 3
$ws.Range("D3").Value = "1.862.38"
$ws.Range("E3").Value = "  -1.00%  "

# Row 4
Set-TextValue "D4" '0.9994'
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
Set-TextValue "D5" '233.76'
$ws.Range("E5").Value = "  -2.18%  "

# Row 6
$ws.Range("E6").Value = "  +0.10%  "

# Row 7
Set-TextValue "D7" '0.4743'
$ws.Range("E7").Value = "  -1.69%  "

# Row 8
Set-TextValue "D8" '0.2744'
$ws.Range("E8").Value = "  -3.39%  "

# Row 9
Set-TextValue "D9" '0.06434'
$ws.Range("E9").Value = "  -1.59%  "

# Row 10
$ws.Range("D10").Value = "1.891.81"
$ws.Range("E10").Value = "  -4.69%  "

# Row 11
Set-TextValue "D11" '0.07448'
$ws.Range("E11").Value = "  -0.52%  "

# Row 12
Set-TextValue "D12" '16.23'
$ws.Range("E12").Value = "  -2.43%  "

# Row 13
Set-TextValue "D13" '5.011'
$ws.Range("E13").Value = "  -2.40%  "

# Row 14
Set-TextValue "D14" '85.47'
$ws.Range("E14").Value = "  -3.86%  "

# Row 15
Set-TextValue "D15" '0.6324'
$ws.Range("E15").Value = "  -5.20%  "

# Row 16
$ws.Range("D16").Value = "30.257.37"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17
Set-TextValue "D17" '1.001'
$ws.Range("E17").Value = "  +0.11%  "

# Row 18
Set-TextValue "D18" '233.20'
$ws.Range("E18").Value = "  -0.98%  "

# Row 19
Set-TextValue "D19" '12.82'
$ws.Range("E19").Value = "  -4.15%  "

# Row 20
Set-TextValue "D20" '0.000007360'
$ws.Range("E20").Value = "  -3.69%  "

# Row 21
$ws.Range("D21").Value = "2.098.14"
$ws.Range("E21").Value = "  -3.08%  "

# Row 22
Set-TextValue "D22" '0.9989'
$ws.Range("E22").Value = "  -0.11%  "

# Row 23
Set-TextValue "D23" '5.088'
$ws.Range("E23").Value = "  -4.65%  "

# Row 24
Set-TextValue "D24" '6.003'
$ws.Range("E24").Value = "  -3.51%  "

# Row 25
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D25" '167.85'
$ws.Range("E25").Value = "  +0.54%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D26" '9.270'
$ws.Range("E26").Value = "  -0.75%  "

# Row 27
Set-TextValue "D27" '17.85'
$ws.Range("E27").Value = "  -4.44%  "

# Row 28
Set-TextValue "D28" '1.876'
$ws.Range("E28").Value = "  -4.27%  "

# Row 29
Set-TextValue "D29" '1.383'
$ws.Range("E29").Value = "  -4.78%  "

# Row 30
Set-TextValue "D30" '0.09984'
$ws.Range("E30").Value = "  +4.05%  "

# Row 31
Set-TextValue "D31" '4.176'
$ws.Range("E31").Value = "  -4.27%  "

# Row 32
Set-TextValue "D32" '3.934'

# Row 33
Set-TextValue "D33" '0.04893'
$ws.Range("E33").Value = "  -3.14%  "

# Row 34
Set-TextValue "D34" '1.145'
$ws.Range("E34").Value = "  -5.95%  "

# Row 35
Set-TextValue "D35" '0.7158'
$ws.Range("E35").Value = "  -5.02%  "

# Row 36
Set-TextValue "D36" '1.0000'
$ws.Range("E36").Value = "  -0.80%  "

# Row 37
Set-TextValue "D37" '2.692'
$ws.Range("E37").Value = "  -0.51%  "

# Row 38
Set-TextValue "D38" '0.01914'
$ws.Range("E38").Value = "  +2.77%  "

# Row 39
Set-TextValue "D39" '2.635'
$ws.Range("E39").Value = "  +0.15%  "

# Row 40
Set-TextValue "D40" '0.9029'
$ws.Range("E40").Value = "  -1.38%  "

# Row 41
Set-TextValue "D41" '1.983'
$ws.Range("E41").Value = "  -5.07%  "

# Row 42
Set-TextValue "D42" '105.84'
$ws.Range("E42").Value = "  -0.60%  "

# Row 43
Set-TextValue "D43" '1.001'
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
Set-TextValue "D44" '0.4108'
$ws.Range("E44").Value = "  -4.32%  "

# Row 45
Set-TextValue "D45" '5.574'
$ws.Range("E45").Value = "  -4.12%  "

# Row 46
Set-TextValue "D46" '7.067'
$ws.Range("E46").Value = "  -5.44%  "

# Row 47
Set-TextValue "D47" '61.23'
$ws.Range("E47").Value = "  -5.61%  "

# Row 48
Set-TextValue "D48" '0.1203'
$ws.Range("E48").Value = "  -6.77%  "

# Row 49
Set-TextValue "D49" '8.705'
$ws.Range("E49").Value = "  -3.47%  "

# Row 50
$ws.Range("E50").Value = "  -6.02%  "

# Row 51
Set-TextValue "D51" '33.12'
$ws.Range("E51").Value = "  -2.50%  "
